$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows after existing row 2 (model_9_6_2), pushing down former
# model_9_6_1 / model_9_6_0 rows, to make room for new models 9_6_3..9_6_8
$ws.Rows.Item(3).EntireRow.Insert()
$ws.Rows.Item(3).EntireRow.Insert()
$ws.Rows.Item(3).EntireRow.Insert()
$ws.Rows.Item(3).EntireRow.Insert()
$ws.Rows.Item(3).EntireRow.Insert()
$ws.Rows.Item(3).EntireRow.Insert()

# Copy column-A formatting (bold, bordered, centered) down onto the new rows
$ws.Range("A2").Copy()
$ws.Range("A3:A8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 2: model_9_6_8
$ws.Range("A2").Value = "model_9_6_8"
$ws.Range("B2").Value = 0.9352443027915913
$ws.Range("C2").Value = 0.8263303822962003
$ws.Range("D2").Value = 0.96064042864996
$ws.Range("E2").Value = 0.8678556736542831
$ws.Range("F2").Value = 0.9692874461994534
$ws.Range("G2").Value = 0.4330216454051976
$ws.Range("H2").Value = 1.161329533260368
$ws.Range("I2").Value = 0.5220328487099741
$ws.Range("J2").Value = 0.1532337483912041
$ws.Range("K2").Value = 0.3376332985505891
$ws.Range("L2").Value = 1.427253679354558
$ws.Range("M2").Value = 0.6580438020414732
$ws.Range("N2").Value = 1.310827346600362
$ws.Range("O2").Value = 0.6860580890328374
$ws.Range("P2").Value = 59.67393512568125
$ws.Range("Q2").Value = 95.02133404685907

# row 3: model_9_6_7
$ws.Range("A3").Value = "model_9_6_7"
$ws.Range("B3").Value = 0.9342127404146074
$ws.Range("C3").Value = 0.8252944876139052
$ws.Range("D3").Value = 0.9607176615396346
$ws.Range("E3").Value = 0.8653096218494853
$ws.Range("F3").Value = 0.9691997551897672
$ws.Range("G3").Value = 0.4399197077699968
$ws.Range("H3").Value = 1.168256565770732
$ws.Range("I3").Value = 0.5210084954452419
$ws.Range("J3").Value = 0.1561861344106133
$ws.Range("K3").Value = 0.3385973149279276
$ws.Range("L3").Value = 1.440360216188453
$ws.Range("M3").Value = 0.6632644327641856
$ws.Range("N3").Value = 1.315778846009884
$ws.Range("O3").Value = 0.691500972813611
$ws.Range("P3").Value = 59.64232610212535
$ws.Range("Q3").Value = 94.98972502330317

# row 4: model_9_6_6
$ws.Range("A4").Value = "model_9_6_6"
$ws.Range("B4").Value = 0.9331589556158159
$ws.Range("C4").Value = 0.8242335427868839
$ws.Range("D4").Value = 0.9607939593395596
$ws.Range("E4").Value = 0.8626916707502391
$ws.Range("F4").Value = 0.969107708074101
$ws.Range("G4").Value = 0.446966371571749
$ws.Range("H4").Value = 1.175351108714225
$ws.Range("I4").Value = 0.519996544438689
$ws.Range("J4").Value = 0.1592218944098195
$ws.Range("K4").Value = 0.3396092194242542
$ws.Range("L4").Value = 1.453579284211246
$ws.Range("M4").Value = 0.6685554364237486
$ws.Range("N4").Value = 1.320837013044084
$ws.Range("O4").Value = 0.6970172254528492
$ws.Range("P4").Value = 59.61054383716341
$ws.Range("Q4").Value = 94.95794275834122

# row 5: model_9_6_5
$ws.Range("A5").Value = "model_9_6_5"
$ws.Range("B5").Value = 0.9320824242600333
$ws.Range("C5").Value = 0.8231471223167731
$ws.Range("D5").Value = 0.96086910751157
$ws.Range("E5").Value = 0.8600016775045866
$ws.Range("F5").Value = 0.9690111679009231
$ws.Range("G5").Value = 0.4541651417048383
$ws.Range("H5").Value = 1.182616007400362
$ws.Range("I5").Value = 0.5189998411473589
$ws.Range("J5").Value = 0.1623411940390747
$ws.Range("K5").Value = 0.3406705175932168
$ws.Range("L5").Value = 1.466913203710486
$ws.Range("M5").Value = 0.6739177558907602
$ws.Range("N5").Value = 1.32600436355184
$ws.Range("O5").Value = 0.702607829961103
$ws.Range("P5").Value = 59.5785887976652
$ws.Range("Q5").Value = 94.92598771884302

# row 6: model_9_6_4
$ws.Range("A6").Value = "model_9_6_4"
$ws.Range("B6").Value = 0.9309824741806386
$ws.Range("C6").Value = 0.8220348433590858
$ws.Range("D6").Value = 0.9609426715873826
$ws.Range("E6").Value = 0.8572381734233737
$ws.Range("F6").Value = 0.9689097951229766
$ws.Range("G6").Value = 0.4615205129505547
$ws.Range("H6").Value = 1.190053821911991
$ws.Range("I6").Value = 0.5180241479997479
$ws.Range("J6").Value = 0.1655457363812919
$ws.Range("K6").Value = 0.3417849421905199
$ws.Range("L6").Value = 1.480365428598592
$ws.Range("M6").Value = 0.6793530105553038
$ws.Range("N6").Value = 1.331284123932935
$ws.Range("O6").Value = 0.7082744746692445
$ws.Range("P6").Value = 59.54645755519359
$ws.Range("Q6").Value = 94.8938564763714

# row 7: model_9_6_3
$ws.Range("A7").Value = "model_9_6_3"
$ws.Range("B7").Value = 0.9298585213484163
$ws.Range("C7").Value = 0.8208962668153986
$ws.Range("D7").Value = 0.961013923717512
$ws.Range("E7").Value = 0.8544010837459085
$ws.Range("F7").Value = 0.9688031467314192
$ws.Range("G7").Value = 0.4690363907149519
$ws.Range("H7").Value = 1.19766748850229
$ws.Range("I7").Value = 0.5170791185903219
$ws.Range("J7").Value = 0.1688356081285104
$ws.Range("K7").Value = 0.3429573633594161
$ws.Range("L7").Value = 1.493937602277077
$ws.Range("M7").Value = 0.6848623151517039
$ws.Range("N7").Value = 1.336679097527602
$ws.Range("O7").Value = 0.7140183217681462
$ws.Range("P7").Value = 59.51414984280643
$ws.Range("Q7").Value = 94.86154876398425

# row 8: model_9_6_2
$ws.Range("A8").Value = "model_9_6_2"
$ws.Range("B8").Value = 0.9287099671834584
$ws.Range("C8").Value = 0.8197310600304539
$ws.Range("D8").Value = 0.9610824290549845
$ws.Range("E8").Value = 0.8514898486591528
$ws.Range("F8").Value = 0.9686909308853908
$ws.Range("G8").Value = 0.4767167777046311
$ws.Range("H8").Value = 1.205459231638499
$ws.Range("I8").Value = 0.5161705203702213
$ws.Range("J8").Value = 0.1722114584364871
$ws.Range("K8").Value = 0.3441909894033541
$ws.Range("L8").Value = 1.507630543753346
$ws.Range("M8").Value = 0.6904467957088591
$ws.Range("N8").Value = 1.3421921575194
$ws.Range("O8").Value = 0.7198405452240878
$ws.Range("P8").Value = 59.48166544374038
$ws.Range("Q8").Value = 94.82906436491818

# row 9: model_9_6_1
$ws.Range("A9").Value = "model_9_6_1"
$ws.Range("B9").Value = 0.9275360895176334
$ws.Range("C9").Value = 0.8185387635378092
$ws.Range("D9").Value = 0.9611477126830839
$ws.Range("E9").Value = 0.8485016228992693
$ws.Range("F9").Value = 0.9685727110360084
$ws.Range("G9").Value = 0.4845665030612124
$ws.Range("H9").Value = 1.213432123774834
$ws.Range("I9").Value = 0.5153046522425495
$ws.Range("J9").Value = 0.1756765866556749
$ws.Range("K9").Value = 0.3454906194491122
$ws.Range("L9").Value = 1.521478403137377
$ws.Range("M9").Value = 0.6961081116185994
$ws.Range("N9").Value = 1.347826770315359
$ws.Range("O9").Value = 0.7257428750726452
$ws.Range("P9").Value = 59.44900119170104
$ws.Range("Q9").Value = 94.79640011287886

# row 10: model_9_6_0
$ws.Range("A10").Value = "model_9_6_0"
$ws.Range("B10").Value = 0.9263362374809045
$ws.Range("C10").Value = 0.8173189706413151
$ws.Range("D10").Value = 0.9612089961690912
$ws.Range("E10").Value = 0.8454358537996909
$ws.Range("F10").Value = 0.9684479884556104
$ws.Range("G10").Value = 0.4925899191556292
$ws.Range("H10").Value = 1.221588884490321
$ws.Range("I10").Value = 0.5144918386960098
$ws.Range("J10").Value = 0.1792316336548271
$ws.Range("K10").Value = 0.3468617361754184
$ws.Range("L10").Value = 1.535846140309702
$ws.Range("M10").Value = 0.70184750420275
$ws.Range("N10").Value = 1.353586060091658
$ws.Range("O10").Value = 0.7317266054812837
$ws.Range("P10").Value = 59.41615651613879
$ws.Range("Q10").Value = 94.7635554373166
